$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Update the password value in B2 from "Admin123" to "admin123"
$ws.Range("B2").Value = "admin123"

# Update the selected cell in the sheet view from B5 to C2
$ws.Range("C2").Select()
